# "Fruta / hortaliza, semanal" weekly refresh.
# The weekly data pull re-ordered the existing records (rows 2-17); the
# per-row field values themselves are unchanged, only which row they land
# on changes. Row 1 (header) and row 18 stay where they are.
#
# Strategy: snapshot the varying columns (D, L, M, N, O, P, Q, S) for every
# data row into memory first, then write them back out in their new order.
# Reading everything up front avoids clobbering a row before it has been
# copied elsewhere (the permutation contains cycles longer than 2).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Maps the NEW row number -> the OLD row number whose data now lives there.
$rowSource = @{
    2  = 16
    3  = 17
    4  = 7
    5  = 8
    6  = 5
    7  = 13
    8  = 15
    9  = 14
    10 = 9
    11 = 10
    12 = 6
    13 = 4
    14 = 3
    15 = 11
    16 = 12
    17 = 2
}

# Snapshot the old values for every data row (2-17) before any writes happen.
# NOTE: use .Value2 (not .Value) for the reads here — in this COM host,
# round-tripping `.Value` through a PowerShell variable captures a broken
# reflection stand-in instead of the real scalar. `.Value2` behaves.
$snapshot = @{}
for ($r = 2; $r -le 17; $r++) {
    $snapshot[$r] = @{
        D = $ws.Cells.Item($r, 4).Value2
        L = $ws.Cells.Item($r, 12).Value2
        M = $ws.Cells.Item($r, 13).Value2
        N = $ws.Cells.Item($r, 14).Value2
        O = $ws.Cells.Item($r, 15).Value2
        P = $ws.Cells.Item($r, 16).Value2
        Q = $ws.Cells.Item($r, 17).Value2
        S = $ws.Cells.Item($r, 19).Value2
    }
}

# Now write each row's new contents from the snapshot taken above.
foreach ($newRow in $rowSource.Keys) {
    $oldRow = $rowSource[$newRow]
    $data = $snapshot[$oldRow]

    $ws.Cells.Item($newRow, 4).Value2 = $data.D
    $ws.Cells.Item($newRow, 12).Value2 = $data.L
    $ws.Cells.Item($newRow, 13).Value2 = $data.M
    $ws.Cells.Item($newRow, 14).Value2 = $data.N
    $ws.Cells.Item($newRow, 15).Value2 = $data.O
    $ws.Cells.Item($newRow, 16).Value2 = $data.P
    $ws.Cells.Item($newRow, 17).Value2 = $data.Q
    $ws.Cells.Item($newRow, 19).Value2 = $data.S
}
